$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header value for new year column K (row 3)
$ws.Range("K3").Value = 2021

# Add the new data values for column K, rows 4-9
$ws.Range("K4").Value = 295
$ws.Range("K5").Value = 163
$ws.Range("K6").Value = 268
$ws.Range("K7").Value = 155
$ws.Range("K8").Value = 27
$ws.Range("K9").Value = 8

# Copy formatting from column J into column K (rows 2-9 only; row 1 has no K cell)
$ws.Range("J2:J9").Copy()
$ws.Range("K2:K9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the active selection to match the recorded cursor position after edit
$ws.Range("L5").Select()
